$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear cells that become empty in the target layout
$ws.Range("A13").Clear()
$ws.Range("B18:C18").Clear()

# Write final cell values (labels + content), matching corrected row alignment
$ws.Range("B1").Value = 'Ementa atual:'
$ws.Range("C1").Value = 'Ementa modificada (dados modificados em vermelho):'
$ws.Range("B2").Value = '8800013'
$ws.Range("C2").Value = '8800013'
$ws.Range("A3").Value = 'Nome:'
$ws.Range("B3").Value = ' Projetos Especiais em Engenharia II'
$ws.Range("C3").Value = ' Projetos Especiais em Engenharia II'
$ws.Range("A4").Value = 'Name:'
$ws.Range("B4").Value = 'Speciais Engeneering Projects II'
$ws.Range("C4").Value = 'Speciais Engeneering Projects II'
$ws.Range("A5").Value = 'Créditos-aula:'
$ws.Range("B5").Value = '2'
$ws.Range("C5").Value = '2'
$ws.Range("A6").Value = 'Créditos-trabalho'
$ws.Range("B6").Value = '2'
$ws.Range("C6").Value = '2'
$ws.Range("A7").Value = 'Carga horária:'
$ws.Range("B7").Value = '90 h'
$ws.Range("C7").Value = '90 h'
$ws.Range("A8").Value = 'Ativação:'
$ws.Range("B8").Value = '01/01/2020'
$ws.Range("C8").Value = '01/01/2020'
$ws.Range("A9").Value = 'Semestre ideal:'
$ws.Range("B9").Value = 'EQD-8'
$ws.Range("C9").Value = 'EQD-8'
$ws.Range("A10").Value = 'Objetivos:'
$ws.Range("B10").Value = 'Levar os alunos a vivenciarem de forma mais aprofundada problemas reais da indústria para, em equipes, apresentarem as possíveis soluções, de forma que, com isso, desenvolvam habilidades transversais fundamentais para sua vida profissional, tais como trabalho em equipe, gerenciamento de projetos, pro atividade, ao mesmo tempo em que consolidam o conhecimento adquirido durante o curso.'
$ws.Range("C10").Value = 'Levar os alunos a vivenciarem de forma mais aprofundada problemas reais da indústria para, em equipes, apresentarem as possíveis soluções, de forma que, com isso, desenvolvam habilidades transversais fundamentais para sua vida profissional, tais como trabalho em equipe, gerenciamento de projetos, pro atividade, ao mesmo tempo em que consolidam o conhecimento adquirido durante o curso.'
$ws.Range("A11").Value = 'Objectives:'
$ws.Range("B11").Value = 'To lead students to experience in-depth real problems of the industry in order to present the possible solutions in teams, so that they develop transversal skills that are fundamental to their professional life, such as teamwork, project management, pro activity, at the same time in which they consolidate the knowledge acquired during the course.'
$ws.Range("C11").Value = 'To lead students to experience in-depth real problems of the industry in order to present the possible solutions in teams, so that they develop transversal skills that are fundamental to their professional life, such as teamwork, project management, pro activity, at the same time in which they consolidate the knowledge acquired during the course.'
$ws.Range("A12").Value = 'Docentes responsáveis:'
$ws.Range("B13").Value = '198273 - Domingos Savio Giordani'
$ws.Range("C13").Value = '198273 - Domingos Savio Giordani'
$ws.Range("A14").Value = 'Programa resumido:'
$ws.Range("B14").Value = 'Formação e trabalho em equipes, Comunicação, Inovação Sistemática, Legislação, Gerenciamento de Projetos. Identificação de Problemas, Formulação do Projeto, Especificação de Problemas, Análise do Conhecimento disponível, Avaliação e Tomada de Decisão, Cronograma, Elaboração de relatórios, Apresentação de Projetos'
$ws.Range("C14").Value = 'Formação e trabalho em equipes, Comunicação, Inovação Sistemática, Legislação, Gerenciamento de Projetos. Identificação de Problemas, Formulação do Projeto, Especificação de Problemas, Análise do Conhecimento disponível, Avaliação e Tomada de Decisão, Cronograma, Elaboração de relatórios, Apresentação de Projetos'
$ws.Range("A15").Value = 'Short syllabus:'
$ws.Range("B15").Value = 'Training and team work, Communication, Systematic Innovation, Legislation, Project Management. Problem Identification, Project Formulation, Problem Specification, Available Knowledge Analysis, Evaluation and Decision Making, Timeline, Reporting, Project Submission'
$ws.Range("C15").Value = 'Training and team work, Communication, Systematic Innovation, Legislation, Project Management. Problem Identification, Project Formulation, Problem Specification, Available Knowledge Analysis, Evaluation and Decision Making, Timeline, Reporting, Project Submission'
$ws.Range("A16").Value = 'Programa:'
$ws.Range("B16").Value = 'Formação e trabalho em equipes e Comunicação – o desenvolvimento das habilidades essenciais para o trabalho em equipes; Inovação Sistemática – desenvolvimento de soluções inovadoras, sistematização e características; Legislação - noções da legislação aplicada à ação empresarial; Gerenciamento de Projetos e Cronograma – Metodologias e esquematizações necessárias com os elementos gerenciais; Identificação de Problemas – sistematização de ações para a localização de causas; Formulação do Projeto – apresentação dos aspectos gerenciais necessários ao desenvolvimento do projeto, Plano de gestão, Estrutura Analítica do Projeto (EAP) etc; Especificação de Problemas – sistematização dos problemas dentro das áreas de conhecimento; Análise do Conhecimento disponível, Avaliação e Tomada de Decisão; Elaboração de relatórios – formatação dentro das normas ABNT; Apresentação de Projetos.'
$ws.Range("C16").Value = 'Formação e trabalho em equipes e Comunicação – o desenvolvimento das habilidades essenciais para o trabalho em equipes; Inovação Sistemática – desenvolvimento de soluções inovadoras, sistematização e características; Legislação - noções da legislação aplicada à ação empresarial; Gerenciamento de Projetos e Cronograma – Metodologias e esquematizações necessárias com os elementos gerenciais; Identificação de Problemas – sistematização de ações para a localização de causas; Formulação do Projeto – apresentação dos aspectos gerenciais necessários ao desenvolvimento do projeto, Plano de gestão, Estrutura Analítica do Projeto (EAP) etc; Especificação de Problemas – sistematização dos problemas dentro das áreas de conhecimento; Análise do Conhecimento disponível, Avaliação e Tomada de Decisão; Elaboração de relatórios – formatação dentro das normas ABNT; Apresentação de Projetos.'
$ws.Range("A17").Value = 'Syllabus:'
$ws.Range("B17").Value = 'Training and work in teams and communication - the development of skills essential to work in teams; Systematic Innovation - development of innovative solutions, systematization and characteristics; Legislation - notions of legislation applied to corporate action; Project Management and Schedule - Methodologies and necessary schematizations with the management elements; Problem Identification - systematization of actions to locate causes; Formulation of the Project - presentation of the managerial aspects necessary for the development of the project, Management Plan, Project Analytical Structure (EAP) etc; Specification of Problems - systematization of problems within the areas of knowledge; Analysis of Available Knowledge, Evaluation and Decision Making; Reporting - formatting within ABNT standards; Presentation of Projects.'
$ws.Range("C17").Value = 'Training and work in teams and communication - the development of skills essential to work in teams; Systematic Innovation - development of innovative solutions, systematization and characteristics; Legislation - notions of legislation applied to corporate action; Project Management and Schedule - Methodologies and necessary schematizations with the management elements; Problem Identification - systematization of actions to locate causes; Formulation of the Project - presentation of the managerial aspects necessary for the development of the project, Management Plan, Project Analytical Structure (EAP) etc; Specification of Problems - systematization of problems within the areas of knowledge; Analysis of Available Knowledge, Evaluation and Decision Making; Reporting - formatting within ABNT standards; Presentation of Projects.'
$ws.Range("A18").Value = 'Avaliação:'
$ws.Range("A19").Value = 'Método:'
$ws.Range("B19").Value = 'Apresentações intermediárias e finais.'
$ws.Range("C19").Value = 'Apresentações intermediárias e finais.'
$ws.Range("A20").Value = 'Critério:'
$ws.Range("B20").Value = 'Serão feitas duas avaliações por uma banca de professores que assistirão às apresentações, as notas serão as médias das notas dadas pelos professores.'
$ws.Range("C20").Value = 'Serão feitas duas avaliações por uma banca de professores que assistirão às apresentações, as notas serão as médias das notas dadas pelos professores.'
$ws.Range("A21").Value = 'Norma de recuperação:'
$ws.Range("B21").Value = 'Reapresentação do último seminário, cuja nota constituirá a nota final da disciplina.'
$ws.Range("C21").Value = 'Reapresentação do último seminário, cuja nota constituirá a nota final da disciplina.'
$ws.Range("A22").Value = 'Bibliografia:'
$ws.Range("B22").Value = 'Gestão de Negócios: Visões e dimensões empresariais da o Organização. Autores: Cruz Jr, J.B., Rocha, J.A.O. e Tachizawa, T.Editora: ATLASGestão Empresarial - de Taylor aos nossos diasAutores: Pereira, M. I. , Autor: Ferreira, A. A. e Reis, A.C. F Editora: THOMSON PIONEIRABaron e Shane: Empreendedorismo: uma visão do processo (EVP), Ed. Thomson, 2006Textos fornecidos pelo professor da disciplinaArtigos extraídos de revistas especializadas na área de gestão e produção.'
$ws.Range("C22").Value = 'Gestão de Negócios: Visões e dimensões empresariais da o Organização. Autores: Cruz Jr, J.B., Rocha, J.A.O. e Tachizawa, T.Editora: ATLASGestão Empresarial - de Taylor aos nossos diasAutores: Pereira, M. I. , Autor: Ferreira, A. A. e Reis, A.C. F Editora: THOMSON PIONEIRABaron e Shane: Empreendedorismo: uma visão do processo (EVP), Ed. Thomson, 2006Textos fornecidos pelo professor da disciplinaArtigos extraídos de revistas especializadas na área de gestão e produção.'

# Row height corrections
$ws.Rows.Item(13).EntireRow.AutoFit()
$ws.Rows.Item(17).RowHeight = 120
$ws.Rows.Item(18).EntireRow.AutoFit()
$ws.Rows.Item(21).RowHeight = 60
$ws.Rows.Item(22).RowHeight = 120
